# Add new columns I (I0) and J (IF) to the worksheet, mirroring the
# style of the existing header row and filling in the per-row values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for the two new columns, styled like the existing header cells.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows 2-29 for columns I (I0) and J (IF).
$values = @{
    2  = @(9, 9)
    3  = @(6, 7)
    4  = @(3, 4)
    5  = @(1, 3)
    6  = @(1, 2)
    7  = @(1, 7)
    8  = @(1, 6)
    9  = @(1, 5)
    10 = @(1, 6)
    11 = @(5, 8)
    12 = @(1, 1)
    13 = @(1, 5)
    14 = @(1, 5)
    15 = @(1, 5)
    16 = @(1, 5)
    17 = @(1, 6)
    18 = @(1, 6)
    19 = @(1, 7)
    20 = @(1, 5)
    21 = @(1, 7)
    22 = @(1, 6)
    23 = @(1, 7)
    24 = @(1, 7)
    25 = @(1, 6)
    26 = @(1, 6)
    27 = @(1, 7)
    28 = @(1, 6)
    29 = @(1, 5)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
